$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 305.66666
$ws.Range("I38").Value = 305.66666
$ws.Range("K38").Value = 916.9999799999999
$ws.Range("M38").Value = -544.9999799999999
$ws.Range("H58").Value = 1287.625
$ws.Range("I58").Value = 292.5
$ws.Range("J58").Value = 1619.3334
$ws.Range("K58").Value = 877.5
$ws.Range("L58").Value = 4858.0002
$ws.Range("M58").Value = -727.5
$ws.Range("N58").Value = -5158.0002
$ws.Range("H62").Value = 4125
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 4350
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 4350
$ws.Range("M62").Value = -3276
$ws.Range("N62").Value = -5598
$ws.Range("H65").Value = 4125
$ws.Range("I65").Value = 3900
$ws.Range("J65").Value = 4350
$ws.Range("K65").Value = 19500
$ws.Range("L65").Value = 21750
$ws.Range("M65").Value = -16380
$ws.Range("N65").Value = -27990
$ws.Range("H92").Value = 1673
$ws.Range("J92").Value = 1498.5
$ws.Range("L92").Value = 1498.5
$ws.Range("N92").Value = -3994.5
$ws.Range("H113").Value = 111152990
$ws.Range("I113").Value = 333342980
$ws.Range("K113").Value = 333342980
$ws.Range("M113").Value = -333339726
$ws.Range("H118").Value = 664.6667
$ws.Range("I118").Value = 664.6667
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1994.0001
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -337.0001
$ws.Range("N118").ClearContents()
$ws.Range("H133").Value = 74966.664
$ws.Range("J133").Value = 74966.664
$ws.Range("L133").Value = 74966.664
$ws.Range("N133").Value = -85086.664
$ws.Range("H138").Value = 4680.268
$ws.Range("I138").Value = 5112.5835
$ws.Range("J138").Value = 4501.3794
$ws.Range("K138").Value = 15337.7505
$ws.Range("L138").Value = 13504.1382
$ws.Range("M138").Value = -10197.7505
$ws.Range("N138").Value = -23784.1382
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 5479.5
$ws.Range("I14").Value = 959
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 959
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -784
$ws.Range("N14").Value = -10350
$ws.Range("H32").Value = 15541.911
$ws.Range("I32").Value = 14485.404
$ws.Range("K32").Value = 14485.404
$ws.Range("M32").Value = -14198.404
$ws.Range("H61").Value = 6596.5864
$ws.Range("I61").Value = 6796.4644
$ws.Range("K61").Value = 6796.4644
$ws.Range("M61").Value = -6584.4644
$ws.Range("H132").Value = 2204.7742
$ws.Range("I132").Value = 1942.5862
$ws.Range("J132").Value = 6006.5
$ws.Range("K132").Value = 5827.7586
$ws.Range("L132").Value = 18019.5
$ws.Range("M132").Value = -3297.7586
$ws.Range("N132").Value = -23079.5
$ws.Range("H136").Value = 6596.5864
$ws.Range("I136").Value = 6796.4644
$ws.Range("K136").Value = 20389.3932
$ws.Range("M136").Value = -17839.3932
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2068.6365
$ws.Range("I94").Value = 2183.6875
$ws.Range("J94").Value = 1761.8334
$ws.Range("K94").Value = 2183.6875
$ws.Range("L94").Value = 1761.8334
$ws.Range("M94").Value = -1732.6875
$ws.Range("N94").Value = -2663.8334
$ws.Range("H134").Value = 991.65625
$ws.Range("I134").Value = 991.65625
$ws.Range("K134").Value = 2974.96875
$ws.Range("M134").Value = -439.96875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 40000
$ws.Range("J41").Value = 45000
$ws.Range("L41").Value = 45000
$ws.Range("N41").Value = -45856
$ws.Range("H51").Value = 19784.75
$ws.Range("I51").Value = 9669.5
$ws.Range("J51").Value = 29900
$ws.Range("K51").Value = 9669.5
$ws.Range("L51").Value = 29900
$ws.Range("M51").Value = -8933.5
$ws.Range("N51").Value = -31372
$ws.Range("H58").Value = 5937.2383
$ws.Range("I58").Value = 2635
$ws.Range("K58").Value = 2635
$ws.Range("M58").Value = -2432
$ws.Range("H61").Value = 19784.75
$ws.Range("I61").Value = 9669.5
$ws.Range("J61").Value = 29900
$ws.Range("K61").Value = 9669.5
$ws.Range("L61").Value = 29900
$ws.Range("M61").Value = -9321.5
$ws.Range("N61").Value = -30596
$ws.Range("H120").Value = 31530.4
$ws.Range("J120").Value = 31530.4
$ws.Range("L120").Value = 31530.4
$ws.Range("N120").Value = -38788.4
$ws.Range("H132").Value = 7571.143
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 10199.6
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 30598.8
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -35658.8
$ws.Range("H136").Value = 5937.2383
$ws.Range("I136").Value = 2635
$ws.Range("K136").Value = 7905
$ws.Range("M136").Value = -5355
$ws.Range("H141").Value = 173835.3
$ws.Range("I141").Value = 31897.75
$ws.Range("J141").Value = 202222.8
$ws.Range("K141").Value = 31897.75
$ws.Range("L141").Value = 202222.8
$ws.Range("M141").Value = -26717.75
$ws.Range("N141").Value = -212582.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 98
$ws.Range("I47").Value = 98
$ws.Range("K47").Value = 294
$ws.Range("M47").Value = 137
$ws.Range("H50").Value = 774.4706
$ws.Range("I50").Value = 983.625
$ws.Range("J50").Value = 588.55554
$ws.Range("K50").Value = 2950.875
$ws.Range("L50").Value = 1765.66662
$ws.Range("M50").Value = -2469.875
$ws.Range("N50").Value = -2727.66662
$ws.Range("H53").Value = 774.4706
$ws.Range("I53").Value = 983.625
$ws.Range("J53").Value = 588.55554
$ws.Range("K53").Value = 2950.875
$ws.Range("L53").Value = 1765.66662
$ws.Range("M53").Value = -2469.875
$ws.Range("N53").Value = -2727.66662
$ws.Range("H68").Value = 83338800
$ws.Range("I68").Value = 166667260
$ws.Range("J68").Value = 10335.333
$ws.Range("K68").Value = 500001780
$ws.Range("L68").Value = 31005.999
$ws.Range("M68").Value = -500000969
$ws.Range("N68").Value = -32627.999
$ws.Range("H71").Value = 83338800
$ws.Range("I71").Value = 166667260
$ws.Range("J71").Value = 10335.333
$ws.Range("K71").Value = 1500005340
$ws.Range("L71").Value = 93017.997
$ws.Range("M71").Value = -1500001284
$ws.Range("N71").Value = -101129.997
$ws.Range("H86").Value = 202
$ws.Range("I86").Value = 202
$ws.Range("K86").Value = 606
$ws.Range("M86").Value = 580
$ws.Range("H89").Value = 202
$ws.Range("I89").Value = 202
$ws.Range("K89").Value = 1818
$ws.Range("M89").Value = 4110
$ws.Range("H123").Value = 15018.25
$ws.Range("I123").Value = 574.5
$ws.Range("K123").Value = 1723.5
$ws.Range("M123").Value = 726.5
$ws.Range("H129").Value = 373204.2
$ws.Range("I129").Value = 67212.07000000001
$ws.Range("K129").Value = 201636.21
$ws.Range("M129").Value = -196636.21
$ws.Range("H132").Value = 1372.5714
$ws.Range("I132").Value = 977
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 8793
$ws.Range("L132").Value = 17100
$ws.Range("M132").Value = -6263
$ws.Range("N132").Value = -22160
$ws.Range("H141").Value = 9411.6
$ws.Range("I141").Value = 9411.6
$ws.Range("K141").Value = 28234.8
$ws.Range("M141").Value = -23054.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 146143820
$ws.Range("I11").Value = 204500880
$ws.Range("J11").Value = 251200
$ws.Range("K11").Value = 204500880
$ws.Range("L11").Value = 251200
$ws.Range("M11").Value = -204500741
$ws.Range("N11").Value = -251478
$ws.Range("H18").Value = 20000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 20000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20586
$ws.Range("M18").ClearContents()
$ws.Range("H113").Value = 4000000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 4536.1665
$ws.Range("I122").Value = 4536.1665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13608.4995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11158.4995
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4984.1177
$ws.Range("J126").Value = 7749.75
$ws.Range("L126").Value = 23249.25
$ws.Range("N126").Value = -28189.25
$ws.Range("H132").Value = 3986.1843
$ws.Range("I132").Value = 3646.9119
$ws.Range("K132").Value = 10940.7357
$ws.Range("M132").Value = -8410.735700000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4770
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H40").Value = 5597.5386
$ws.Range("I40").Value = 5206.1055
$ws.Range("J40").Value = 6660
$ws.Range("K40").Value = 5206.1055
$ws.Range("L40").Value = 6660
$ws.Range("M40").Value = -5070.1055
$ws.Range("N40").Value = -6932
$ws.Range("H122").Value = 3865.3103
$ws.Range("I122").Value = 3904.3333
$ws.Range("J122").Value = 3678
$ws.Range("K122").Value = 11712.9999
$ws.Range("L122").Value = 11034
$ws.Range("M122").Value = -9262.999899999999
$ws.Range("N122").Value = -15934
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3790.875
$ws.Range("I122").Value = 2757.5715
$ws.Range("J122").Value = 5237.5
$ws.Range("K122").Value = 8272.7145
$ws.Range("L122").Value = 15712.5
$ws.Range("M122").Value = -5822.7145
$ws.Range("N122").Value = -20612.5
$ws.Range("H132").Value = 3080.0195
$ws.Range("I132").Value = 2683.3777
$ws.Range("K132").Value = 8050.1331
$ws.Range("M132").Value = -5520.1331
$ws.Range("H135").Value = 129499
$ws.Range("J135").Value = 129499
$ws.Range("L135").Value = 129499
$ws.Range("N135").Value = -139639
